$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2515.65
$ws.Range("I17").Value = 1099
$ws.Range("K17").Value = 3297
$ws.Range("M17").Value = -3129

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 235489.67
$ws.Range("J70").Value = 235489.67
$ws.Range("L70").Value = 706469.01
$ws.Range("N70").Value = -707009.01

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 235489.67
$ws.Range("J73").Value = 235489.67
$ws.Range("L73").Value = 706469.01
$ws.Range("N73").Value = -708341.01

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 634.94116
$ws.Range("I92").Value = 581.6
$ws.Range("K92").Value = 581.6
$ws.Range("M92").Value = 666.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2215.3333
$ws.Range("I98").Value = 1973.8334
$ws.Range("K98").Value = 1973.8334
$ws.Range("M98").Value = -475.8334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1324.0869
$ws.Range("I100").Value = 953
$ws.Range("J100").Value = 2660
$ws.Range("K100").Value = 953
$ws.Range("L100").Value = 2660
$ws.Range("M100").Value = -412
$ws.Range("N100").Value = -3742

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 18956.54
$ws.Range("I106").Value = 18956.54
$ws.Range("K106").Value = 18956.54
$ws.Range("M106").Value = -18325.54

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2215.3333
$ws.Range("I122").Value = 1973.8334
$ws.Range("K122").Value = 5921.5002
$ws.Range("M122").Value = -3471.5002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 38999.5
$ws.Range("I137").Value = 38999.5
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 116998.5
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -114448.5
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 12277.474
$ws.Range("I138").Value = 11047.5
$ws.Range("J138").Value = 12716.75
$ws.Range("K138").Value = 33142.5
$ws.Range("L138").Value = 38150.25
$ws.Range("M138").Value = -28002.5
$ws.Range("N138").Value = -48430.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 22484
$ws.Range("J24").Value = 22484
$ws.Range("L24").Value = 22484
$ws.Range("N24").Value = -23232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 23877.777
$ws.Range("I32").Value = 17979.316
$ws.Range("J32").Value = 30470.176
$ws.Range("K32").Value = 17979.316
$ws.Range("L32").Value = 30470.176
$ws.Range("M32").Value = -17692.316
$ws.Range("N32").Value = -31044.176

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1413
$ws.Range("I61").Value = 884
$ws.Range("K61").Value = 884
$ws.Range("M61").Value = -672

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 22484
$ws.Range("J100").Value = 22484
$ws.Range("L100").Value = 22484
$ws.Range("N100").Value = -24648

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3902.4666
$ws.Range("I122").Value = 3008.5
$ws.Range("J122").Value = 4924.143
$ws.Range("K122").Value = 9025.5
$ws.Range("L122").Value = 14772.429
$ws.Range("M122").Value = -6575.5
$ws.Range("N122").Value = -19672.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1933.1052
$ws.Range("I132").Value = 1839.375
$ws.Range("K132").Value = 5518.125
$ws.Range("M132").Value = -2988.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1413
$ws.Range("I136").Value = 884
$ws.Range("K136").Value = 2652
$ws.Range("M136").Value = -102

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 674.9286
$ws.Range("I22").Value = 675
$ws.Range("J22").Value = 674.75
$ws.Range("K22").Value = 675
$ws.Range("L22").Value = 674.75
$ws.Range("M22").Value = -502
$ws.Range("N22").Value = -1020.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7953.7144
$ws.Range("I86").Value = 7999.6665
$ws.Range("K86").Value = 7999.6665
$ws.Range("M86").Value = -6876.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 7953.7144
$ws.Range("I89").Value = 7999.6665
$ws.Range("K89").Value = 39998.3325
$ws.Range("M89").Value = -34382.3325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 151.66667
$ws.Range("I94").Value = 151.66667
$ws.Range("K94").Value = 151.66667
$ws.Range("M94").Value = 299.33333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3770.6667
$ws.Range("I134").Value = 2241.3333
$ws.Range("K134").Value = 6723.999899999999
$ws.Range("M134").Value = -4188.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 30000
$ws.Range("J20").Value = 30000
$ws.Range("L20").Value = 30000
$ws.Range("N20").Value = -30472

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 21249.5
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 21249.5
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 21249.5
$ws.Range("M28").ClearContents()
$ws.Range("N28").Value = -21739.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H30").Value = 30000
$ws.Range("J30").Value = 30000
$ws.Range("L30").Value = 30000
$ws.Range("N30").Value = -30182

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3783.75
$ws.Range("J58").Value = 6676.5713
$ws.Range("L58").Value = 6676.5713
$ws.Range("N58").Value = -7082.5713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I94").Value = 800
$ws.Range("J94").Value = 1158
$ws.Range("K94").Value = 800
$ws.Range("L94").Value = 1158
$ws.Range("M94").Value = -349
$ws.Range("N94").Value = -2060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H128").Value = 30000
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2358.2812
$ws.Range("I132").Value = 2211.7036
$ws.Range("K132").Value = 6635.110799999999
$ws.Range("M132").Value = -4105.110799999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4641.1
$ws.Range("I134").Value = 4375
$ws.Range("J134").Value = 4907.2
$ws.Range("K134").Value = 13125
$ws.Range("L134").Value = 14721.6
$ws.Range("M134").Value = -10590
$ws.Range("N134").Value = -19791.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 3783.75
$ws.Range("J136").Value = 6676.5713
$ws.Range("L136").Value = 20029.7139
$ws.Range("N136").Value = -25129.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 3750.5
$ws.Range("I34").Value = 3334
$ws.Range("K34").Value = 10002
$ws.Range("M34").Value = -9918

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 1019
$ws.Range("J52").Value = 1019
$ws.Range("L52").Value = 3057
$ws.Range("N52").Value = -3589

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 101347.4
$ws.Range("J55").Value = 1857.1428
$ws.Range("L55").Value = 5571.428400000001
$ws.Range("N55").Value = -5925.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 87.875
$ws.Range("I97").Value = 82.25
$ws.Range("J97").Value = 93.5
$ws.Range("K97").Value = 246.75
$ws.Range("L97").Value = 280.5
$ws.Range("M97").Value = 249.25
$ws.Range("N97").Value = -1272.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 746.2308
$ws.Range("J114").Value = 765.1
$ws.Range("L114").Value = 2295.3
$ws.Range("N114").Value = -8803.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3586.3
$ws.Range("I131").Value = 1470
$ws.Range("J131").Value = 4115.375
$ws.Range("K131").Value = 4410
$ws.Range("L131").Value = 12346.125
$ws.Range("M131").Value = 630
$ws.Range("N131").Value = -22426.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 5223.2
$ws.Range("I140").Value = 4939.4443
$ws.Range("K140").Value = 14818.3329
$ws.Range("M140").Value = -9638.332900000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5163.7144
$ws.Range("J43").Value = 9882.666999999999
$ws.Range("L43").Value = 9882.666999999999
$ws.Range("N43").Value = -10184.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 231
$ws.Range("I107").Value = 203.6
$ws.Range("J107").Value = 299.5
$ws.Range("K107").Value = 203.6
$ws.Range("L107").Value = 299.5
$ws.Range("M107").Value = 1716.4
$ws.Range("N107").Value = -4139.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4640.65
$ws.Range("I132").Value = 3924.1538
$ws.Range("K132").Value = 11772.4614
$ws.Range("M132").Value = -9242.4614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2592.25
$ws.Range("I46").Value = 1335.2
$ws.Range("K46").Value = 1335.2
$ws.Range("M46").Value = -1147.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1742
$ws.Range("I93").Value = 1594.1666
$ws.Range("K93").Value = 1594.1666
$ws.Range("M93").Value = -346.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4300.364
$ws.Range("I132").Value = 4300.364
$ws.Range("K132").Value = 12901.092
$ws.Range("M132").Value = -10371.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6394.65
$ws.Range("I62").Value = 5332.3335
$ws.Range("J62").Value = 6582.1177
$ws.Range("K62").Value = 5332.3335
$ws.Range("L62").Value = 6582.1177
$ws.Range("M62").Value = -4708.3335
$ws.Range("N62").Value = -7830.1177

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 6394.65
$ws.Range("I65").Value = 5332.3335
$ws.Range("J65").Value = 6582.1177
$ws.Range("K65").Value = 26661.6675
$ws.Range("L65").Value = 32910.5885
$ws.Range("M65").Value = -23541.6675
$ws.Range("N65").Value = -39150.5885

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2676.5557
$ws.Range("I113").Value = 2027.6
$ws.Range("J113").Value = 3487.75
$ws.Range("K113").Value = 6082.799999999999
$ws.Range("L113").Value = 10463.25
$ws.Range("M113").Value = -3912.799999999999
$ws.Range("N113").Value = -14803.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2990.75
$ws.Range("I132").Value = 2990.75
$ws.Range("K132").Value = 8972.25
$ws.Range("M132").Value = -6442.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 66793.625
$ws.Range("I136").Value = 3070.1
$ws.Range("K136").Value = 9210.299999999999
$ws.Range("M136").Value = -6660.299999999999
